$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-14 08:16:58"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-14 08:16:47"
$wsZhCn.Range("K2").Value = "2016-10-14 08:17:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-14 08:16:58"
$wsDeDe.Range("K2").Value = "2016-10-14 08:17:41"
